$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new tracking row for the "Non-Overlapping Intervals" question,
# solved using a greedy approach.
$ws.Rows("12").Insert()

$ws.Range("A12").Value = 46076
$ws.Range("B12").Value = "Non-Overlapping Intervals"
$ws.Hyperlinks.Add($ws.Range("C12"), "https://leetcode.com/problems/non-overlapping-intervals/") | Out-Null

# Match the date-number-format cell style already used elsewhere in column C/A.
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

$wb.Save()
